$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (45181 -> 45182) for every data row (rows 2 through 44).
for ($row = 2; $row -le 44; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value = 45182
    }
}
